# Mifos Automation Excel edit:
#   - "Repayment schedule" sheet: insert a new (blank) column before the
#     existing "Late"/Outstanding columns (i.e. before column N), copying
#     the column-width formatting of the column immediately to its left
#     (column M) onto the freshly inserted column.
#   - Make "Repayment schedule" the active sheet/tab (it was previously
#     "Floating Interest Rates"), and leave the last used cell selected
#     on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position N (column 14). Excel shifts the
# existing N:P columns (Late / Outstanding / Disbursement-like "Original"
# amount) one slot to the right, to O:Q.
$ws.Columns.Item(14).Insert()

# The newly inserted column picks up the default width; match it to the
# width of the column directly to its left (M, "In Advance") like Excel
# does when you insert via the column header context menu.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Switch the active sheet to "Repayment schedule" and select the last
# cell the user left the cursor on.
$ws.Activate()
$ws.Range("R8").Select()
